# Update test case to include a tie-break.
# - Row 5 (Ballot #4) changes 1st/2nd/3rd choice from A,B,C to B,C,D
# - Row 8 (Ballot #7) changes 1st/2nd/3rd choice from B,C,D to C,D,undervote
# - The last data row (Ballot #10 / row 11) is removed entirely, shrinking
#   the table to 9 ballots (rows 2-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5: Ballot #4 now ranks B, C, D
$ws.Range("B5").Value = "B"
$ws.Range("C5").Value = "C"
$ws.Range("D5").Value = "D"

# Update row 8: Ballot #7 now ranks C, D, undervote
$ws.Range("B8").Value = "C"
$ws.Range("C8").Value = "D"
$ws.Range("D8").Value = "undervote"

# Remove the last data row (row 11, Ballot #10) entirely
$ws.Rows("11").Delete()

# Update the active selection to match the saved state
$ws.Range("C11").Select()
